# Update NATMI ligand-receptor TPM-derived statistics with newly
# recomputed values (commit: "update scripts wuth new tpm").
#
# Only the derived expression / specificity columns (G:T, excluding the
# untouched K and L columns) change; the grouping columns (A:F) and the
# detection-rate / cell-count columns that did not change stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (ECs -> ECs)
$ws.Range("G2").Value = 0.8901166666666667
$ws.Range("H2").Value = 2.67035
$ws.Range("I2").Value = 0.001454906651255627
$ws.Range("J2").Value = 0.001454906651255627
$ws.Range("M2").Value = 1.639301666666666
$ws.Range("N2").Value = 4.917904999999999
$ws.Range("O2").Value = 0.1482499788910018
$ws.Range("P2").Value = 0.1482499788910018
$ws.Range("Q2").Value = 1.459169735194444
$ws.Range("R2").Value = 13.13252761675
$ws.Range("S2").Value = 0.0002156898803370248
$ws.Range("T2").Value = 0.0002156898803370249

# Row 3 (ECs -> FAPs)
$ws.Range("G3").Value = 0.8901166666666667
$ws.Range("H3").Value = 2.67035
$ws.Range("I3").Value = 0.001454906651255627
$ws.Range("J3").Value = 0.001454906651255627
$ws.Range("O3").Value = 0.7250300748586421
$ws.Range("P3").Value = 0.7250300748586421
$ws.Range("Q3").Value = 7.136202988044446
$ws.Range("R3").Value = 64.22582689240001
$ws.Range("S3").Value = 0.001054851078272204
$ws.Range("T3").Value = 0.001054851078272204

# Row 4 (ECs -> MuSCs)
$ws.Range("G4").Value = 0.8901166666666667
$ws.Range("H4").Value = 2.67035
$ws.Range("I4").Value = 0.001454906651255627
$ws.Range("J4").Value = 0.001454906651255627
$ws.Range("M4").Value = 1.401229333333333
$ws.Range("N4").Value = 4.203688
$ws.Range("O4").Value = 0.1267199462503561
$ws.Range("P4").Value = 0.1267199462503561
$ws.Range("Q4").Value = 1.247257583422222
$ws.Range("R4").Value = 11.2253182508
$ws.Range("S4").Value = 0.0001843656926463987
$ws.Range("T4").Value = 0.0001843656926463987

# Row 5 (FAPs -> ECs)
$ws.Range("I5").Value = 0.9436845966890257
$ws.Range("J5").Value = 0.9436845966890258
$ws.Range("M5").Value = 1.639301666666666
$ws.Range("N5").Value = 4.917904999999999
$ws.Range("O5").Value = 0.1482499788910018
$ws.Range("P5").Value = 0.1482499788910018
$ws.Range("Q5").Value = 946.4497271143915
$ws.Range("R5").Value = 8518.047544029525
$ws.Range("S5").Value = 0.1399012215389116
$ws.Range("T5").Value = 0.1399012215389116

# Row 6 (FAPs -> FAPs)
$ws.Range("I6").Value = 0.9436845966890257
$ws.Range("J6").Value = 0.9436845966890258
$ws.Range("O6").Value = 0.7250300748586421
$ws.Range("P6").Value = 0.7250300748586421
$ws.Range("S6").Value = 0.6841997137803918
$ws.Range("T6").Value = 0.6841997137803919

# Row 7 (FAPs -> MuSCs)
$ws.Range("I7").Value = 0.9436845966890257
$ws.Range("J7").Value = 0.9436845966890258
$ws.Range("M7").Value = 1.401229333333333
$ws.Range("N7").Value = 4.203688
$ws.Range("O7").Value = 0.1267199462503561
$ws.Range("P7").Value = 0.1267199462503561
$ws.Range("Q7").Value = 808.9988237824932
$ws.Range("R7").Value = 7280.989414042439
$ws.Range("S7").Value = 0.1195836613697223
$ws.Range("T7").Value = 0.1195836613697223

# Row 8 (MuSCs -> ECs)
$ws.Range("G8").Value = 33.563832
$ws.Range("H8").Value = 100.691496
$ws.Range("I8").Value = 0.05486049665971852
$ws.Range("J8").Value = 0.05486049665971854
$ws.Range("M8").Value = 1.639301666666666
$ws.Range("N8").Value = 4.917904999999999
$ws.Range("O8").Value = 0.1482499788910018
$ws.Range("P8").Value = 0.1482499788910018
$ws.Range("Q8").Value = 55.02124573731999
$ws.Range("R8").Value = 495.1912116358799
$ws.Range("S8").Value = 0.008133067471753144
$ws.Range("T8").Value = 0.008133067471753148

# Row 9 (MuSCs -> FAPs)
$ws.Range("G9").Value = 33.563832
$ws.Range("H9").Value = 100.691496
$ws.Range("I9").Value = 0.05486049665971852
$ws.Range("J9").Value = 0.05486049665971854
$ws.Range("O9").Value = 0.7250300748586421
$ws.Range("P9").Value = 0.7250300748586421
$ws.Range("Q9").Value = 269.086432350016
$ws.Range("R9").Value = 2421.777891150144
$ws.Range("S9").Value = 0.03977550999997801
$ws.Range("T9").Value = 0.03977550999997802

# Row 10 (MuSCs -> MuSCs)
$ws.Range("G10").Value = 33.563832
$ws.Range("H10").Value = 100.691496
$ws.Range("I10").Value = 0.05486049665971852
$ws.Range("J10").Value = 0.05486049665971854
$ws.Range("M10").Value = 1.401229333333333
$ws.Range("N10").Value = 4.203688
$ws.Range("O10").Value = 0.1267199462503561
$ws.Range("P10").Value = 0.1267199462503561
$ws.Range("Q10").Value = 47.030625937472
$ws.Range("R10").Value = 423.275633437248
$ws.Range("S10").Value = 0.006951919187987371
$ws.Range("T10").Value = 0.006951919187987375
